# feat: Make sidebar non-fixed and allow full-width content when collapsed
# (Data update applied to the "Business Contact Information" sheet of cpms_data.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business Contact Information")

# --- Row 2 fix-up: I2/J2 were stored as numbers; re-enter them as text    ---
# (the displayed value is unchanged: "21321" / "213123", only the cell's
# underlying data type switches from numeric to text)
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "21321"
$ws.Range("I2").Style = "Normal"

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "213123"
$ws.Range("J2").Style = "Normal"

# --- Row 3: new data row ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Registered"
$ws.Range("C3").Value = "REGION XI (DAVAO REGION)"
$ws.Range("D3").Value = "DAVAO DEL SUR"
$ws.Range("E3").Value = "DAVAO CITY"
$ws.Range("F3").Value = "Acacia"
$ws.Range("G3").Value = "21344tytr3"
$ws.Range("H3").Value = "1435y643524"
$ws.Range("I3").Value = "4324577i76543"
$ws.Range("J3").Value = "esadfdgghhjgfsd"
